# Update layout of zip files: add "Get your Box" / C11 item as a new row
# above the existing "Order internal swag" (C10) row on the Planning sheet,
# pushing the remaining Content/communications-plan rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 17 (shifts old rows 17-44 down to 18-45,
# inheriting formatting/height from the row above, same as Excel's
# "Insert Sheet Rows").
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row with the new task.
$ws.Range("A17").Value2 = "C5"
$ws.Range("B17").Value2 = "Get your Box"
$ws.Range("C17").Value2 = "Send an email to derek.teay@ibm.com with the number of developers that will be attending your Call for Code day and your shipping address"

# Match the row height used for the new row in the authored workbook.
$ws.Rows.Item(17).RowHeight = 33

# Reflect the author's final selection (the whole new row selected).
$ws.Rows.Item(17).Select()
